$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("events")

# B2 keeps its "date" styling, but the date format changes from the
# built-in m/d/yyyy format to the custom "d mmm" format.
$ws.Range("B2").NumberFormat = "d mmm"

# B3 and B4 lose the date formatting entirely - copy the plain
# (no-number-format) style already used by C2 in the same row onto them.
$ws.Range("C2").Copy()
$ws.Range("B3").PasteSpecial(-4122)
$ws.Range("B4").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# All three event dates move from 1 Nov 2025 (45962) to 26 Nov 2025 (45987).
$ws.Range("B2").Value = 45987
$ws.Range("B3").Value = 45987
$ws.Range("B4").Value = 45987
